# Update the date line and every multiplication problem in the table
# to match the new "output generated" snapshot.

$d = $word.ActiveDocument

$replacements = @(
    @("2024-05-20 Monday", "2024-05-21 Tuesday"),
    @("431×8=", "387×4="),
    @("153×3=", "905×7="),
    @("315×6=", "861×4="),
    @("716×2=", "895×5="),
    @("443×6=", "504×3="),
    @("379×7=", "396×8="),
    @("932×9=", "230×8="),
    @("982×6=", "622×9="),
    @("685×8=", "236×7="),
    @("465×7=", "259×5="),
    @("616×4=", "662×3="),
    @("754×3=", "219×5="),
    @("862×2=", "218×4="),
    @("333×8=", "638×4="),
    @("944×3=", "313×3="),
    @("222×9=", "491×7="),
    @("990×2=", "780×7="),
    @("807×4=", "244×2="),
    @("190×5=", "274×7="),
    @("757×5=", "696×7="),
    @("166×2=", "884×7="),
    @("722×6=", "595×7="),
    @("642×6=", "765×8="),
    @("932×2=", "709×3="),
    @("390×7=", "496×9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
